# Updates for Version 2.0 in documentation
#
# - Replace the "=General!B4" version-lookup formula in the three document
#   sheets with a literal "2.0" text value (same visual style, no formula).
# - Bump the "Date" column from 2018-05-10 to 2018-05-21.
# - Move the active sheet / selection from "Project notebook" to
#   "User documentation", and update the remembered selections on the
#   sheets that were visited along the way.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Project notebook", "Developer handbook", "User documentation")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $c2 = $ws.Range("C2")
    # Build the literal "2.0" as a formula result first, then collapse it to
    # a plain value via copy / paste-values so the cell keeps its existing
    # style (s="3") instead of picking up a new number-formatted style the
    # way a direct Value="2.0" assignment would (Excel would infer the
    # number 2 from "2.0" typed directly into a General-formatted cell).
    $c2.Formula = "=""2.0"""
    $c2.Copy()
    $c2.PasteSpecial(-4163)  # xlPasteValues

    # New document date: 2018-05-10 -> 2018-05-21
    $ws.Range("D2").Value = "2018-05-21"
}

$excel.CutCopyMode = 0

# Walk the sheets in the same order the saved selections appear in the
# target workbook, ending on "User documentation" so it becomes the
# active/selected tab.
$wsProjectNotebook = $wb.Worksheets.Item("Project notebook")
$wsProjectNotebook.Activate()
$wsProjectNotebook.Range("D2").Select()

$wsDeveloperHandbook = $wb.Worksheets.Item("Developer handbook")
$wsDeveloperHandbook.Activate()
$wsDeveloperHandbook.Range("C2:D2").Select()

$wsUserDocumentation = $wb.Worksheets.Item("User documentation")
$wsUserDocumentation.Activate()
$wsUserDocumentation.Range("G6").Select()
